$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.288329
$ws.Range("H2").Value = 117.864987
$ws.Range("I2").Value = 0.632237668435316
$ws.Range("J2").Value = 0.632237668435316
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 1128.524772173026
$ws.Range("R2").Value = 10156.72294955724
$ws.Range("S2").Value = 0.267642586830699
$ws.Range("T2").Value = 0.267642586830699

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.288329
$ws.Range("H3").Value = 117.864987
$ws.Range("I3").Value = 0.632237668435316
$ws.Range("J3").Value = 0.632237668435316
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 1201.038808836708
$ws.Range("R3").Value = 10809.34927953037
$ws.Range("S3").Value = 0.2848401219072513
$ws.Range("T3").Value = 0.2848401219072513

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.288329
$ws.Range("H4").Value = 117.864987
$ws.Range("I4").Value = 0.632237668435316
$ws.Range("J4").Value = 0.632237668435316
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 336.2897093020279
$ws.Range("R4").Value = 3026.607383718252
$ws.Range("S4").Value = 0.07975495969736572
$ws.Range("T4").Value = 0.07975495969736572

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.344283
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3112930657211948
$ws.Range("J5").Value = 0.3112930657211947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 555.6485379010533
$ws.Range("R5").Value = 5000.83684110948
$ws.Range("S5").Value = 0.1317784205712876
$ws.Range("T5").Value = 0.1317784205712876

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.344283
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3112930657211948
$ws.Range("J6").Value = 0.3112930657211947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 591.3520682470404
$ws.Range("R6").Value = 5322.168614223363
$ws.Range("S6").Value = 0.1402459220886785
$ws.Range("T6").Value = 0.1402459220886785

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.344283
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3112930657211948
$ws.Range("J7").Value = 0.3112930657211947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 165.578009355556
$ws.Range("R7").Value = 1490.202084200004
$ws.Range("S7").Value = 0.03926872306122861
$ws.Range("T7").Value = 0.0392687230612286

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.509096666666667
$ws.Range("H8").Value = 10.52729
$ws.Range("I8").Value = 0.05646926584348937
$ws.Range("J8").Value = 0.05646926584348937
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 100.7959008967555
$ws.Range("R8").Value = 907.1631080707999
$ws.Range("S8").Value = 0.0239049033952462
$ws.Range("T8").Value = 0.02390490339524621

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.509096666666667
$ws.Range("H9").Value = 10.52729
$ws.Range("I9").Value = 0.05646926584348937
$ws.Range("J9").Value = 0.05646926584348937
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 107.2726020143589
$ws.Range("R9").Value = 965.45341812923
$ws.Range("S9").Value = 0.0254409273124765
$ws.Range("T9").Value = 0.0254409273124765

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.509096666666667
$ws.Range("H10").Value = 10.52729
$ws.Range("I10").Value = 0.05646926584348937
$ws.Range("J10").Value = 0.05646926584348937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 30.03622520942666
$ws.Range("R10").Value = 270.32602688484
$ws.Range("S10").Value = 0.007123435135766664
$ws.Range("T10").Value = 0.007123435135766664
